$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.549.83'
$ws.Range('E2').Value = '  +0.49%  '

$ws.Range('D3').Value = '1.737.74'
$ws.Range('E3').Value = '  +0.51%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9996'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.13%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.90'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.49%  '

$ws.Range('E6').Value = '  -0.12%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4930'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +2.94%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2670'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.43%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06299'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.26%  '

$ws.Range('D10').Value = '1.732.77'
$ws.Range('E10').Value = '  +0.07%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07046'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -1.16%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.72'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.36%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.593'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.94%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6129'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.38%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '77.49'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.80%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.000'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.08%  '

$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000007354'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +6.43%  '

$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '26.537.96'
$ws.Range('E18').Value = '  +0.34%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.000'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.11%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.55'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.26%  '

$ws.Range('D21').Value = '1.953.80'
$ws.Range('E21').Value = '  -0.12%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.588'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.37%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.717'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -1.93%  '

$ws.Range('E24').Value = '  -1.31%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '140.07'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +2.75%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.47'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.67%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.420'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.27%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '107.94'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +1.28%  '

$ws.Range('E29').Value = '  -1.42%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.041'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.68%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08063'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.59%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.719'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.09%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04600'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +1.46%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.612'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.25%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.012'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +2.59%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6368'
$ws.Range('D36').ClearFormats()

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.8970'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -3.68%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.012'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.92%  '

$ws.Range('E39').Value = '  +0.14%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.003'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.30%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.01507'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +0.22%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '101.96'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -6.88%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.401'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -4.22%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3905'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.26%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '6.869'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.88%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1188'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.12%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05399'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +1.26%  '

$ws.Range('B48').Value = 'Elrond'
$ws.Range('C48').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '30.53'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.83%  '

$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.771'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -1.61%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.270'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.01%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '51.79'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.99%  '
